$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.159980058670044
$ws.Range("B1").Value = 2.407708168029785
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.383955478668213
$ws.Range("E1").Value = 1.228777170181274
